$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CCSaleData")

# Update Date (column B) for rows 2-115 with new timestamps,
# and flip Result (column A) from Fail to Pass where applicable.
$dates = @(
    "Sun Oct 31 13:20:26 EDT 2021"
    "Sun Oct 31 13:21:34 EDT 2021"
    "Sun Oct 31 13:22:35 EDT 2021"
    "Sun Oct 31 13:23:37 EDT 2021"
    "Sun Oct 31 13:24:37 EDT 2021"
    "Sun Oct 31 13:25:38 EDT 2021"
    "Sun Oct 31 13:26:39 EDT 2021"
    "Sun Oct 31 13:27:40 EDT 2021"
    "Sun Oct 31 13:28:40 EDT 2021"
    "Sun Oct 31 13:29:42 EDT 2021"
    "Sun Oct 31 13:30:43 EDT 2021"
    "Sun Oct 31 13:31:43 EDT 2021"
    "Sun Oct 31 13:32:45 EDT 2021"
    "Sun Oct 31 13:33:45 EDT 2021"
    "Sun Oct 31 13:34:45 EDT 2021"
    "Sun Oct 31 13:35:47 EDT 2021"
    "Sun Oct 31 13:36:47 EDT 2021"
    "Sun Oct 31 13:37:48 EDT 2021"
    "Sun Oct 31 13:38:48 EDT 2021"
    "Sun Oct 31 13:39:49 EDT 2021"
    "Sun Oct 31 13:40:50 EDT 2021"
    "Sun Oct 31 13:41:51 EDT 2021"
    "Sun Oct 31 13:42:51 EDT 2021"
    "Sun Oct 31 13:43:52 EDT 2021"
    "Sun Oct 31 13:44:55 EDT 2021"
    "Sun Oct 31 13:45:55 EDT 2021"
    "Sun Oct 31 13:46:56 EDT 2021"
    "Sun Oct 31 13:47:58 EDT 2021"
    "Sun Oct 31 13:48:59 EDT 2021"
    "Sun Oct 31 13:50:02 EDT 2021"
    "Sun Oct 31 13:51:04 EDT 2021"
    "Sun Oct 31 13:52:06 EDT 2021"
    "Sun Oct 31 13:53:06 EDT 2021"
    "Sun Oct 31 13:54:09 EDT 2021"
    "Sun Oct 31 13:55:09 EDT 2021"
    "Sun Oct 31 13:56:10 EDT 2021"
    "Sun Oct 31 13:57:13 EDT 2021"
    "Sun Oct 31 13:58:17 EDT 2021"
    "Sun Oct 31 13:59:18 EDT 2021"
    "Sun Oct 31 14:00:22 EDT 2021"
    "Sun Oct 31 14:01:26 EDT 2021"
    "Sun Oct 31 14:02:27 EDT 2021"
    "Sun Oct 31 14:03:30 EDT 2021"
    "Sun Oct 31 14:04:31 EDT 2021"
    "Sun Oct 31 14:05:34 EDT 2021"
    "Sun Oct 31 14:06:36 EDT 2021"
    "Sun Oct 31 14:07:37 EDT 2021"
    "Sun Oct 31 14:08:40 EDT 2021"
    "Sun Oct 31 14:09:43 EDT 2021"
    "Sun Oct 31 14:10:43 EDT 2021"
    "Sun Oct 31 14:11:45 EDT 2021"
    "Sun Oct 31 14:12:46 EDT 2021"
    "Sun Oct 31 14:13:48 EDT 2021"
    "Sun Oct 31 14:14:48 EDT 2021"
    "Sun Oct 31 14:15:50 EDT 2021"
    "Sun Oct 31 14:16:52 EDT 2021"
    "Sun Oct 31 14:17:54 EDT 2021"
    "Sun Oct 31 14:18:55 EDT 2021"
    "Sun Oct 31 14:19:58 EDT 2021"
    "Sun Oct 31 14:20:59 EDT 2021"
    "Sun Oct 31 14:22:01 EDT 2021"
    "Sun Oct 31 14:23:02 EDT 2021"
    "Sun Oct 31 14:24:04 EDT 2021"
    "Sun Oct 31 14:25:06 EDT 2021"
    "Sun Oct 31 14:26:07 EDT 2021"
    "Sun Oct 31 14:27:09 EDT 2021"
    "Sun Oct 31 14:28:12 EDT 2021"
    "Sun Oct 31 14:29:13 EDT 2021"
    "Sun Oct 31 14:30:16 EDT 2021"
    "Sun Oct 31 14:31:18 EDT 2021"
    "Sun Oct 31 14:32:20 EDT 2021"
    "Sun Oct 31 14:33:22 EDT 2021"
    "Sun Oct 31 14:34:23 EDT 2021"
    "Sun Oct 31 14:35:24 EDT 2021"
    "Sun Oct 31 14:36:25 EDT 2021"
    "Sun Oct 31 14:37:26 EDT 2021"
    "Sun Oct 31 14:38:27 EDT 2021"
    "Sun Oct 31 14:39:28 EDT 2021"
    "Sun Oct 31 14:40:30 EDT 2021"
    "Sun Oct 31 14:41:31 EDT 2021"
    "Sun Oct 31 14:42:33 EDT 2021"
    "Sun Oct 31 14:43:34 EDT 2021"
    "Sun Oct 31 14:44:36 EDT 2021"
    "Sun Oct 31 14:45:37 EDT 2021"
    "Sun Oct 31 14:46:39 EDT 2021"
    "Sun Oct 31 14:47:41 EDT 2021"
    "Sun Oct 31 14:48:43 EDT 2021"
    "Sun Oct 31 14:49:45 EDT 2021"
    "Sun Oct 31 14:50:47 EDT 2021"
    "Sun Oct 31 14:51:49 EDT 2021"
    "Sun Oct 31 14:52:50 EDT 2021"
    "Sun Oct 31 14:53:52 EDT 2021"
    "Sun Oct 31 14:54:54 EDT 2021"
    "Sun Oct 31 14:55:56 EDT 2021"
    "Sun Oct 31 14:56:58 EDT 2021"
    "Sun Oct 31 14:58:00 EDT 2021"
    "Sun Oct 31 14:59:02 EDT 2021"
    "Sun Oct 31 15:00:03 EDT 2021"
    "Sun Oct 31 15:01:04 EDT 2021"
    "Sun Oct 31 15:02:05 EDT 2021"
    "Sun Oct 31 15:03:08 EDT 2021"
    "Sun Oct 31 15:04:10 EDT 2021"
    "Sun Oct 31 15:05:12 EDT 2021"
    "Sun Oct 31 15:06:12 EDT 2021"
    "Sun Oct 31 15:07:12 EDT 2021"
    "Sun Oct 31 15:08:13 EDT 2021"
    "Sun Oct 31 15:09:14 EDT 2021"
    "Sun Oct 31 15:10:15 EDT 2021"
    "Sun Oct 31 15:11:14 EDT 2021"
    "Sun Oct 31 15:12:15 EDT 2021"
    "Sun Oct 31 15:13:17 EDT 2021"
    "Sun Oct 31 15:14:19 EDT 2021"
    "Sun Oct 31 15:15:19 EDT 2021"
    "Sun Oct 31 15:16:21 EDT 2021"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $dates[$i]
}

$passRows = @(40,42,44,46,48,64,66,68,70,72,88,90,92,94,96,101,102,103)
foreach ($r in $passRows) {
    $ws.Cells.Item($r, 1).Value = "Pass"
}
